$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data per diff
$ws.Range("D2").Value = '42.656.01'
$ws.Range("E2").Value = '  -0.39%  '

$ws.Range("D3").Value = '2.545.56'
$ws.Range("E3").Value = '  +0.10%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.73'
$ws.Range("E5").Value = '  +3.31%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '95.07'
$ws.Range("E6").Value = '  -2.70%  '

$ws.Range("E7").Value = '  +0.43%  '

$ws.Range("E8").Value = '  -0.04%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.536'
$ws.Range("E9").Value = '  -1.48%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.36'
$ws.Range("E10").Value = '  -1.18%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0816'
$ws.Range("E11").Value = '  -1.51%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.60'
$ws.Range("E12").Value = '  +0.15%  '

$ws.Range("E13").Value = '  -0.50%  '

$ws.Range("D14").Value = '2.928.96'
$ws.Range("E14").Value = '  -0.22%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.67'
$ws.Range("E15").Value = '  +4.04%  '

$ws.Range("D16").Value = '2.526.82'
$ws.Range("E16").Value = '  -2.34%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.867'
$ws.Range("E17").Value = '  +0.06%  '

$ws.Range("D18").Value = '42.678.37'
$ws.Range("E18").Value = '  -0.44%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.10'
$ws.Range("E19").Value = '  -1.65%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.68'
$ws.Range("E20").Value = '  +1.56%  '

$ws.Range("D21").Value = '0.0₃0971'
$ws.Range("E21").Value = '  -1.84%  '

$ws.Range("E22").Value = '  -1.00%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '254.97'
$ws.Range("E23").Value = '  -0.41%  '

$ws.Range("E24").Value = '  +0.41%  '

$ws.Range("E25").Value = '  -1.49%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '27.64'
$ws.Range("E26").Value = '  -1.58%  '

$ws.Range("E27").Value = '  +0.21%  '

$ws.Range("E28").Value = '  +2.70%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '39.52'
$ws.Range("E29").Value = '  +4.36%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '10.09'
$ws.Range("E30").Value = '  -0.79%  '

$ws.Range("E31").Value = '  -1.65%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '155.59'
$ws.Range("E32").Value = '  -1.48%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.43'
$ws.Range("E33").Value = '  +3.36%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '19.50'
$ws.Range("E34").Value = '  -0.06%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.15'
$ws.Range("E35").Value = '  +0.80%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0792'
$ws.Range("E36").Value = '  -0.67%  '

$ws.Range("E37").Value = '  -0.42%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.112'
$ws.Range("E38").Value = '  -3.84%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '24.62'
$ws.Range("E39").Value = '  -4.31%  '

$ws.Range("E41").Value = '  +7.13%  '

$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.85'
$ws.Range("E42").Value = '  -0.73%  '

$ws.Range("B43").Value = 'NEARProtocol'
$ws.Range("C43").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.38'
$ws.Range("E43").Value = '  -0.69%  '

$ws.Range("E44").Value = '  -0.43%  '

$ws.Range("E45").Value = '  -0.02%  '

$ws.Range("D46").Value = '2.050.55'
$ws.Range("E46").Value = '  -2.00%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '85.27'
$ws.Range("E47").Value = '  -2.62%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.94'
$ws.Range("E48").Value = '  +0.19%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '75.27'
$ws.Range("E49").Value = '  +0.70%  '

$ws.Range("D50").Value = '2.785.80'
$ws.Range("E50").Value = '  -0.25%  '

$ws.Range("E51").Value = '  -0.35%  '
